$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5..66 down to 6..67
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44882
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 100112032
$ws.Range("G5").Value = "Zapallo italiano"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 6500
$ws.Range("N5").Value = "`$/caja 50 unidades"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 130
$ws.Range("Q5").Value = 50
$ws.Range("R5").Value = "Hortaliza"
